# Adiciona a coluna "Confirmar Senha" logo apos a coluna "Senha" na planilha
# "Contas", replicando o valor de Senha, atualiza o nome de usuario da
# primeira conta de teste e ajusta a formatacao/largura da coluna Senha.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contas")

# Insere uma nova coluna antes da antiga coluna D ("Primeiro Nome"), logo
# depois da coluna "Senha" (C) -- isso desloca as colunas D em diante uma
# posicao para a direita automaticamente.
$ws.Columns.Item(4).Insert()

# Preenche o cabecalho e os valores da nova coluna "Confirmar Senha" com o
# mesmo conteudo da coluna "Senha".
$ws.Range("D1").Value = "Confirmar Senha"
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = $ws.Range("C3").Value2

# Atualiza o usuario de teste da primeira conta da massa de dados.
$ws.Range("A2").Value = "AmimGers"

# Aplica o mesmo estilo usado nas celulas de Email (estilo Hiperlink) para as
# celulas de Senha, igual ao que foi feito na coluna nova.
$ws.Range("C2").Style = "Hiperlink"
$ws.Range("C3").Style = "Hiperlink"

# Ajusta a largura da coluna "Senha" (mais proximo possivel do autofit).
$ws.Columns.Item(3).ColumnWidth = 8.666666666666666
